$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 is untouched and keeps the workbook-default (unstyled) format;
# we borrow its Style object below to strip the transient "quote-prefix"
# formatting Excel applies when a numeric-looking literal is forced to text.
$plainStyle = $ws.Range("A1").Style

$ws.Range('D2').Value = '63.721.34'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '3.395.62'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'568.31"
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').Value = "'161.46"
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.401.97'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = "'0.549"
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = '  -5.96%  '
$ws.Range('D10').Value = "'7.26"
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D12').Value = "'0.421"
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = '  -5.09%  '
$ws.Range('D13').Value = '3.992.87'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = "'26.81"
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('E16').Value = '  -4.38%  '
$ws.Range('D17').Value = '63.836.13'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '3.393.47'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = "'6.08"
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('D20').Value = "'13.52"
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').Value = "'374.82"
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').Value = "'7.72"
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = "'70.90"
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('D25').Value = "'0.513"
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = '  -6.30%  '
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('D27').Value = "'9.43"
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  -4.02%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = "'6.06"
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').Value = "'1.38"
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  -6.28%  '
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').Value = "'22.81"
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = "'7.04"
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('E35').Value = '  -4.71%  '
$ws.Range('D36').Value = "'159.79"
$ws.Range('D36').Style = $plainStyle
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('D37').Value = "'0.859"
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  +10.36%  '
$ws.Range('E38').Value = '  -4.73%  '
$ws.Range('D39').Value = "'0.0723"
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = "'25.73"
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.759.02'
$ws.Range('E41').Value = '  -4.57%  '
$ws.Range('D42').Value = "'42.67"
$ws.Range('D42').Style = $plainStyle
$ws.Range('E43').Value = '  -3.25%  '
$ws.Range('D44').Value = "'26.02"
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').Value = "'4.37"
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  -3.69%  '
$ws.Range('D46').Value = "'0.0304"
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('E47').Value = '  +5.98%  '
$ws.Range('D48').Value = "'328.53"
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  +4.00%  '
$ws.Range('E49').Value = '  -4.08%  '
$ws.Range('D50').Value = "'6.27"
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = '  -3.47%  '
$ws.Range('E51').Value = '  -3.17%  '
